$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 33333614
$ws.Range("I11").Value = 33333614
$ws.Range("K11").Value = 33333614
$ws.Range("M11").Value = -33333474

$ws.Range("H98").Value = 1524.3684
$ws.Range("I98").Value = 1586.0588
$ws.Range("K98").Value = 1586.0588
$ws.Range("M98").Value = -88.05880000000002

$ws.Range("H106").Value = 2858.4285
$ws.Range("I106").Value = 2701.5
$ws.Range("J106").Value = 3800
$ws.Range("K106").Value = 2701.5
$ws.Range("L106").Value = 3800
$ws.Range("M106").Value = -2070.5
$ws.Range("N106").Value = -5062

$ws.Range("H122").Value = 1524.3684
$ws.Range("I122").Value = 1586.0588
$ws.Range("K122").Value = 4758.1764
$ws.Range("M122").Value = -2308.1764


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3211.25
$ws.Range("I32").Value = 1716.6072
$ws.Range("J32").Value = 10186.25
$ws.Range("K32").Value = 1716.6072
$ws.Range("L32").Value = 10186.25
$ws.Range("M32").Value = -1429.6072
$ws.Range("N32").Value = -10760.25

$ws.Range("H74").Value = 853.9677
$ws.Range("I74").Value = 795.4074000000001
$ws.Range("J74").Value = 1249.25
$ws.Range("K74").Value = 795.4074000000001
$ws.Range("L74").Value = 1249.25
$ws.Range("M74").Value = 78.59259999999995
$ws.Range("N74").Value = -2997.25

$ws.Range("H77").Value = 853.9677
$ws.Range("I77").Value = 795.4074000000001
$ws.Range("J77").Value = 1249.25
$ws.Range("K77").Value = 3977.037
$ws.Range("L77").Value = 6246.25
$ws.Range("M77").Value = 390.9629999999997
$ws.Range("N77").Value = -14982.25

$ws.Range("H88").Value = 3213.6924
$ws.Range("I88").Value = 2259.8572
$ws.Range("J88").Value = 4326.5
$ws.Range("K88").Value = 2259.8572
$ws.Range("L88").Value = 4326.5
$ws.Range("M88").Value = -1853.8572
$ws.Range("N88").Value = -5138.5

$ws.Range("H91").Value = 3213.6924
$ws.Range("I91").Value = 2259.8572
$ws.Range("J91").Value = 4326.5
$ws.Range("K91").Value = 2259.8572
$ws.Range("L91").Value = 4326.5
$ws.Range("M91").Value = -855.8571999999999
$ws.Range("N91").Value = -7134.5

$ws.Range("H109").Value = 45937
$ws.Range("J109").Value = 45937
$ws.Range("L109").Value = 45937
$ws.Range("N109").Value = -48711

$ws.Range("H130").Value = 90000
$ws.Range("J130").Value = 90000
$ws.Range("L130").Value = 90000
$ws.Range("N130").Value = -100040

$ws.Range("H132").Value = 1406.7084
$ws.Range("I132").Value = 1080.1364
$ws.Range("K132").Value = 3240.4092
$ws.Range("M132").Value = -710.4092000000001


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 84888.125
$ws.Range("I86").Value = 1487.8667
$ws.Range("J86").Value = 223888.56
$ws.Range("K86").Value = 1487.8667
$ws.Range("L86").Value = 223888.56
$ws.Range("M86").Value = -364.8667
$ws.Range("N86").Value = -226134.56

$ws.Range("H89").Value = 84888.125
$ws.Range("I89").Value = 1487.8667
$ws.Range("J89").Value = 223888.56
$ws.Range("K89").Value = 7439.333500000001
$ws.Range("L89").Value = 1119442.8
$ws.Range("M89").Value = -1823.333500000001
$ws.Range("N89").Value = -1130674.8

$ws.Range("H134").Value = 9065.612999999999
$ws.Range("J134").Value = 8909
$ws.Range("L134").Value = 26727
$ws.Range("N134").Value = -31797


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2728.0571
$ws.Range("I31").Value = 1925.2916
$ws.Range("J31").Value = 4479.5454
$ws.Range("K31").Value = 1925.2916
$ws.Range("L31").Value = 4479.5454
$ws.Range("M31").Value = -1630.2916
$ws.Range("N31").Value = -5069.5454

$ws.Range("H34").Value = 2728.0571
$ws.Range("I34").Value = 1925.2916
$ws.Range("J34").Value = 4479.5454
$ws.Range("K34").Value = 1925.2916
$ws.Range("L34").Value = 4479.5454
$ws.Range("M34").Value = -1723.2916
$ws.Range("N34").Value = -4883.5454

$ws.Range("H74").Value = 26190
$ws.Range("J74").Value = 26190
$ws.Range("L74").Value = 26190
$ws.Range("N74").Value = -27938

$ws.Range("H77").Value = 26190
$ws.Range("J77").Value = 26190
$ws.Range("L77").Value = 78570
$ws.Range("N77").Value = -87306

$ws.Range("H99").Value = 1668833
$ws.Range("I99").Value = 3334666
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 3334666
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -3333168
$ws.Range("N99").Value = -5996

$ws.Range("H126").Value = 1668833
$ws.Range("I126").Value = 3334666
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 10003998
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -10001528
$ws.Range("N126").Value = -13940

$ws.Range("H134").Value = 926.6667
$ws.Range("I134").Value = 926.6667
$ws.Range("K134").Value = 2780.0001
$ws.Range("M134").Value = -245.0001000000002

$ws.Range("H141").Value = 48741.5
$ws.Range("J141").Value = 48741.5
$ws.Range("L141").Value = 48741.5
$ws.Range("N141").Value = -59101.5


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 710.7143
$ws.Range("J122").Value = 904.2222
$ws.Range("L122").Value = 8137.999800000001
$ws.Range("N122").Value = -13037.9998

$ws.Range("H129").Value = 61463
$ws.Range("J129").Value = 91850
$ws.Range("L129").Value = 275550
$ws.Range("N129").Value = -285550


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2538.45
$ws.Range("I102").Value = 2529.5833
$ws.Range("K102").Value = 2529.5833
$ws.Range("M102").Value = -907.5832999999998


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9428.286
$ws.Range("I40").Value = 1999.6666
$ws.Range("K40").Value = 1999.6666
$ws.Range("M40").Value = -1863.6666

$ws.Range("H106").Value = 15874.5
$ws.Range("J106").Value = 15874.5
$ws.Range("L106").Value = 15874.5
$ws.Range("N106").Value = -18398.5

$ws.Range("H132").Value = 1835.9678
$ws.Range("I132").Value = 1653.875
$ws.Range("K132").Value = 4961.625
$ws.Range("M132").Value = -2431.625

$ws.Range("H136").Value = 3007.111
$ws.Range("I136").Value = 1694
$ws.Range("J136").Value = 5633.3335
$ws.Range("K136").Value = 5082
$ws.Range("L136").Value = 16900.0005
$ws.Range("M136").Value = -2532
$ws.Range("N136").Value = -22000.0005


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 43913.332
$ws.Range("J95").Value = 43913.332
$ws.Range("L95").Value = 43913.332
$ws.Range("N95").Value = -49405.332

